$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "53.965.12"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.245.86"
$ws.Range("E3").Value = "  +2.86%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "493.48"
$ws.Range("E5").Value = "  +2.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.43"
$ws.Range("E6").Value = "  +3.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.526"
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0952"
$ws.Range("E9").Value = "  +4.57%  "
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("E11").Value = "  +3.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.64"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.648.35"
$ws.Range("E13").Value = "  +3.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.67"
$ws.Range("E14").Value = "  +3.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "53.900.24"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.263.12"
$ws.Range("E17").Value = "  +3.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.00"
$ws.Range("E18").Value = "  +5.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.08"
$ws.Range("E19").Value = "  +4.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.41"
$ws.Range("E20").Value = "  +6.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "298.97"
$ws.Range("E21").Value = "  +2.50%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("E23").Value = "  -3.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.80"
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.369"
$ws.Range("E26").Value = "  +2.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.357.98"
$ws.Range("E27").Value = "  +3.22%  "
$ws.Range("E28").Value = "  +3.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.02"
$ws.Range("E29").Value = "  +1.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.16"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("E31").Value = "  +2.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0682"
$ws.Range("E32").Value = "  +4.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.83"
$ws.Range("E33").Value = "  +3.49%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.63"
$ws.Range("E37").Value = "  +2.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.901"
$ws.Range("E38").Value = "  +11.15%  "
$ws.Range("E39").Value = "  +3.61%  "
$ws.Range("E40").Value = "  +4.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.65"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  +3.68%  "
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.34"
$ws.Range("E44").Value = "  +3.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.92"
$ws.Range("E45").Value = "  +4.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "125.10"
$ws.Range("E46").Value = "  +2.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0886"
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.540"
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "236.38"
$ws.Range("E49").Value = "  +4.16%  "
$ws.Range("E50").Value = "  +3.57%  "
$ws.Range("E51").Value = "  +2.04%  "
